$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 13.33
$ws.Range("E4").Value = 12.919
$ws.Range("E7").Value = 13.35
$ws.Range("E8").Value = 12.913
$ws.Range("D11").Value = -8.342000000000002
$ws.Range("D12").Value = -8.038
$ws.Range("E12").Value = 13.338
$ws.Range("E14").Value = 12.911
$ws.Range("D15").Value = -7.897
$ws.Range("E22").Value = 12.862
